# Add the "Chemistry-Syllabus" worksheet (new tab at the end of the workbook)
# and populate it with the syllabus table, matching the committed change.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Chemistry-Syllabus"

# Header row (SlNo / Topics)
$ws.Range("A1").Value = "SlNo"
$ws.Range("B1").Value = "Topics"

$topics = @(
    'Some Basic Concepts of Chemistry',
    'Solid State',
    'Solutions',
    'Electrochemistry',
    'Chemical Kinetics',
    'Surface Chemistry',
    'Structure of Atom',
    'Classification of elements and periodicity of properties',
    'Chemical Bonding and Molecular State',
    'States of Matter Gases and Liquids',
    'Thermodynamics',
    'Equilibrium',
    'Redox Reactions',
    'Hydrogen',
    'S-Block of elements (Alkali and Alkaline Earth Materials)',
    'Some P-Block Elements',
    'Environmental Chemistry',
    'P-Block Elements',
    'D and F Block Elements',
    'Coordination Compounds ',
    'Haloalkanes and Haloarenes ',
    'Alcohols, Phenols and Ethers ',
    'Organic compounds containing Nitrogen  ',
    'Biomolecules ',
    'Polymers  ',
    'Chemistry in Everyday Life'
)

# Data starts on row 3 (row 2 is left blank) and runs through row 28.
# The first three topic rows are filled in before the two name columns
# (C1/D1) are added, matching the order the workbook was authored in.
for ($i = 0; $i -lt 3; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $topics[$i]
}

$ws.Range("C1").Value = "Sahesta"
$ws.Range("D1").Value = "Prottoy"

for ($i = 3; $i -lt $topics.Length; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $topics[$i]
}

# Column A (the SlNo column) is centred, both header and data rows.
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A3:A28").HorizontalAlignment = -4108

# Column widths, to match the authored sheet as closely as possible.
$ws.Columns("A:A").ColumnWidth = 4.77734375
$ws.Columns("B:B").ColumnWidth = 46.33203125

# Scroll / selection state left by the author.
$ws.Range("B25").Select()

$wb.Save()
